# Add Event Module TestCases:
#   Insert a new "Events" worksheet between "Employee" and "Reward",
#   populate it with the event test data, size its columns, and make
#   it the active/selected sheet (at cell C3).

$wb = $excel.ActiveWorkbook

$empSheet = $wb.Worksheets.Item("Employee")

# Insert the new sheet right after "Employee" (i.e. before "Reward").
$eventsSheet = $wb.Worksheets.Add($null, $empSheet)
$eventsSheet.Name = "Events"

# Column sizing (approx. characters; engine stores widths on a 1/6-char grid).
$eventsSheet.Columns.Item(1).ColumnWidth = 11.5
$eventsSheet.Columns.Item(2).ColumnWidth = 93.33333333333333
$eventsSheet.Columns.Item(3).ColumnWidth = 34.333333333333336

# Header row.
$eventsSheet.Range("A1").Value = "eventname"
$eventsSheet.Range("B1").Value = "eventdescription"
$eventsSheet.Range("C1").Value = "eventdocument"

# Row 2: Code Ninjas event.
$eventsSheet.Range("A2").Value = "code ninjas"
# Row 3: Hacker Throne event.
$eventsSheet.Range("A3").Value = "hacker throne"

$eventsSheet.Range("B2").Value = "Have a blast building awesome video games and developing ninja coding skills on the path from white to black belt."
$eventsSheet.Range("B3").Value = "A hackathon is a fast-paced collaborative event lasting 24-48 hours where people work on engineering projects."

$eventsSheet.Range("C2").Value = "C:/Users/mygoa/Pictures/codeNinjas.png"
$eventsSheet.Range("C3").Value = "C:/Users/mygoa/Pictures/codeNinjas.png"

# Make Events the active sheet with C3 selected (matches tabSelected + selection).
$eventsSheet.Activate()
$eventsSheet.Range("C3").Select()
